$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" - the report generator re-ran and the row
# that used to belong to e7c36651-26d0-49ee-b1cf-7f40f5dd3b38 and the row
# that used to belong to 6e6f5a0a-f847-4138-9cfe-7cdb61058920 swapped places
# (rows 8 and 9) on every sheet, and e7c36651's handoff timestamp moved
# forward since it was (re-)handed off after 6e6f5a0a.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (Path And Name, hyperlinked),
# E/F (zh-cn / de-de status) and G (Latest HO Xliff Generate Date).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A8").Value = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
$wsOverview.Range("E8").Value = "Ready for handoff"
$wsOverview.Range("F8").Value = "Ready for handoff"
$wsOverview.Range("G8").Value = "2016-12-16 08:12:14"

$wsOverview.Range("A9").Value = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-12-16 08:24:01"

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$8') {
        $h.TextToDisplay = "e2e\6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
    } elseif ($h.Range.Address() -eq '$B$9') {
        $h.TextToDisplay = "e2e\e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
    }
}

# ---------------------------------------------------------------------------
# Sheets "zh-cn" / "de-de": column A (Source File Name, hyperlinked),
# C (Status), G (Latest Handoff File) and H (Latest Handoff Datetime).
# The two locale sheets carry different xlf/date values, so they are
# parameterised below.
# ---------------------------------------------------------------------------
$localeSheets = @(
    @{ Name = "zh-cn"; Row8G = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.fa7b938b8ca0282e071b9dfae621037cafe4c44e.zh-cn.xlf"; Row8H = "2016-12-16 08:12:01"; Row9G = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.ebe39ec3a11a72ce0470bafa1cc822a30a67b978.zh-cn.xlf"; Row9H = "2016-12-16 08:23:46" },
    @{ Name = "de-de"; Row8G = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.fa7b938b8ca0282e071b9dfae621037cafe4c44e.de-de.xlf"; Row8H = "2016-12-16 08:12:14"; Row9G = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.ebe39ec3a11a72ce0470bafa1cc822a30a67b978.de-de.xlf"; Row9H = "2016-12-16 08:24:01" }
)

foreach ($info in $localeSheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    $ws.Range("A8").Value = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
    $ws.Range("C8").Value = "Ready for handoff"
    $ws.Range("G8").Value = $info.Row8G
    $ws.Range("H8").Value = $info.Row8H

    $ws.Range("A9").Value = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
    $ws.Range("C9").Value = "Ready for handoff"
    $ws.Range("G9").Value = $info.Row9G
    $ws.Range("H9").Value = $info.Row9H

    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$A$8') {
            $h.TextToDisplay = "6e6f5a0a-f847-4138-9cfe-7cdb61058920.md"
        } elseif ($h.Range.Address() -eq '$A$9') {
            $h.TextToDisplay = "e7c36651-26d0-49ee-b1cf-7f40f5dd3b38.md"
        }
    }
}
